# Merge the split "<id>...</id>" runs (e.g. "<id>" + "p054v_1" + "</id>")
# into a single run per paragraph, keeping the formatting (Courier New,
# color 7f6000, size 9pt) of the surrounding tag runs.
#
# Only the plain "<id>p054v_N</id>" paragraphs are affected; the similarly
# shaped "<id>fig_p054v_1</id>" paragraph (different inner-run formatting)
# is left untouched, matching the source diff.

$d = $word.ActiveDocument

$ids = @("p054v_1", "p054v_2", "p054v_3", "p054v_4", "p054v_5")

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    # Paragraph text includes a trailing paragraph-mark (CR) character.
    $t = $p.Range.Text.TrimEnd([char]13)

    foreach ($id in $ids) {
        $target = "<id>" + $id + "</id>"
        if ($t -eq $target) {
            $full = $p.Range.Duplicate
            $full.Start = $p.Range.Start
            $full.End = $p.Range.End
            $full.Text = $target
        }
    }
}
